# formInfo workbook rework:
#  - rename "Case2" sheet to "TestMessage"
#  - Settings!B2 now points at the relative lib path instead of the old
#    absolute Windows path
#  - TestMessage (ex-Case2) drops its "arg_booFormOnTop" test row and
#    uses longer sample message/title text for the remaining rows
#  - Case1 becomes the active sheet/tab

$wb = $excel.ActiveWorkbook

$wsSettings = $wb.Worksheets.Item("Settings")
$wsCase1    = $wb.Worksheets.Item("Case1")
$wsCase2    = $wb.Worksheets.Item("Case2")

# Rename the third sheet.
$wsCase2.Name = "TestMessage"

# Settings sheet: shorten the formInfo.xaml path.
$wsSettings.Range("B2").Value = "lib\formInfo.xaml"

# TestMessage sheet: remove the arg_booFormOnTop row (row 3), which shifts
# arg_strRobotName / arg_strInitMessage / arg_strMessageTitle up one row.
$wsCase2.Rows.Item(3).Delete()

# Update the sample values for the (now shifted) Init Message / Message
# Title rows to the longer test strings.
$wsCase2.Range("D4").Value = "Test Message longer here"
$wsCase2.Range("D5").Value = "Test Message Title longer here"

# Update selections on the non-active sheets.
$wsSettings.Range("B2").Select()
$wsCase2.Range("A4").Select()

# Make Case1 the active sheet/tab with A4 selected.
$wsCase1.Activate()
$wsCase1.Range("A4").Select()
